$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 110: id=109, x_pos=7, y_pos=6.5, type=intersection
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = 7
$ws.Range("C110").Value = 6.5
$ws.Range("D110").Value = "intersection"

# Row 111: id=110, x_pos=6.5, y_pos=7, type=intersection
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = 6.5
$ws.Range("C111").Value = 7
$ws.Range("D111").Value = "intersection"

# Row 112: id=111, x_pos=4, y_pos=7.5, type=intersection
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = 4
$ws.Range("C112").Value = 7.5
$ws.Range("D112").Value = "intersection"

# Row 113: id=112, x_pos=7, y_pos=7, type=D (new departure depot)
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = 7
$ws.Range("C113").Value = 7
$ws.Range("D113").Value = "D"

# Row 114: id=113, x_pos=4, y_pos=8, type=A (new arrival depot)
$ws.Range("A114").Value = 113
$ws.Range("B114").Value = 4
$ws.Range("C114").Value = 8
$ws.Range("D114").Value = "A"

# Move the view/selection down to the newly added data, as in the saved file
$null = $ws.Range("G110").Select()
